# Update countries & provincias Spain
#
# Applies the day-over-day data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp in A1.
#  - Refreshes case/recovery/death counts for several countries whose
#    rows did not otherwise move (B-H numeric columns).
#  - Three country names newly enter the "top" listing ahead of their
#    former neighbours (Grecia before Croacia, Trinidad y Tobago before
#    Sri Lanka), which shifts the country label (column A) for the
#    rows that follow while carrying each country's own updated figures
#    along with it.
#  - Two existing rows (Santa Lucia/Timor Oriental and
#    Montserrat/Islas Malvinas) swap with their neighbour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 16 de Septiembre de 2020 a las 18:26'
$ws.Cells.Item(4, 2).Value = 6798195
$ws.Cells.Item(4, 3).Value = 10048
$ws.Cells.Item(4, 4).Value = 4070239
$ws.Cells.Item(4, 5).Value = 2527277
$ws.Cells.Item(4, 7).Value = 482
$ws.Cells.Item(4, 8).Value = 200679
$ws.Cells.Item(5, 2).Value = 5079933
$ws.Cells.Item(5, 3).Value = 61899
$ws.Cells.Item(5, 4).Value = 3995071
$ws.Cells.Item(5, 5).Value = 1002193
$ws.Cells.Item(5, 7).Value = 578
$ws.Cells.Item(5, 8).Value = 82669
$ws.Cells.Item(17, 2).Value = 378219
$ws.Cells.Item(17, 3).Value = 3991
$ws.Cells.Item(17, 7).Value = 20
$ws.Cells.Item(17, 8).Value = 41684
$ws.Cells.Item(23, 2).Value = 291442
$ws.Cells.Item(23, 3).Value = 1452
$ws.Cells.Item(23, 4).Value = 215265
$ws.Cells.Item(23, 5).Value = 40532
$ws.Cells.Item(23, 7).Value = 12
$ws.Cells.Item(23, 8).Value = 35645
$ws.Cells.Item(25, 2).Value = 266541
$ws.Cells.Item(25, 3).Value = 1697
$ws.Cells.Item(25, 5).Value = 17994
$ws.Cells.Item(29, 2).Value = 139421
$ws.Cells.Item(29, 3).Value = 618
$ws.Cells.Item(29, 4).Value = 121994
$ws.Cells.Item(29, 5).Value = 8234
$ws.Cells.Item(29, 7).Value = 5
$ws.Cells.Item(29, 8).Value = 9193
$ws.Cells.Item(57, 2).Value = 57514
$ws.Cells.Item(57, 3).Value = 26
$ws.Cells.Item(57, 5).Value = 603
$ws.Cells.Item(60, 2).Value = 48966
$ws.Cells.Item(60, 3).Value = 232
$ws.Cells.Item(60, 4).Value = 34517
$ws.Cells.Item(60, 5).Value = 12804
$ws.Cells.Item(60, 7).Value = 13
$ws.Cells.Item(60, 8).Value = 1645
$ws.Cells.Item(66, 2).Value = 40186
$ws.Cells.Item(66, 3).Value = 1290
$ws.Cells.Item(66, 4).Value = 22896
$ws.Cells.Item(66, 5).Value = 16809
$ws.Cells.Item(66, 7).Value = 5
$ws.Cells.Item(66, 8).Value = 481
$ws.Cells.Item(89, 1).Value = 'Grecia'
$ws.Cells.Item(89, 2).Value = 14041
$ws.Cells.Item(89, 3).Value = 311
$ws.Cells.Item(89, 4).Value = 3804
$ws.Cells.Item(89, 5).Value = 9921
$ws.Cells.Item(89, 7).Value = 3
$ws.Cells.Item(89, 8).Value = 316
$ws.Cells.Item(90, 1).Value = 'Croacia'
$ws.Cells.Item(90, 2).Value = 14029
$ws.Cells.Item(90, 3).Value = 280
$ws.Cells.Item(90, 4).Value = 11690
$ws.Cells.Item(90, 5).Value = 2103
$ws.Cells.Item(90, 7).Value = 6
$ws.Cells.Item(90, 8).Value = 236
$ws.Cells.Item(91, 1).Value = 'Zambia'
$ws.Cells.Item(91, 2).Value = 13887
$ws.Cells.Item(91, 3).Value = 68
$ws.Cells.Item(91, 4).Value = 12869
$ws.Cells.Item(91, 5).Value = 692
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = 326
$ws.Cells.Item(125, 2).Value = 4164
$ws.Cells.Item(125, 3).Value = 122
$ws.Cells.Item(125, 4).Value = 1180
$ws.Cells.Item(125, 5).Value = 2938
$ws.Cells.Item(136, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(136, 2).Value = 3293
$ws.Cells.Item(136, 3).Value = 70
$ws.Cells.Item(136, 4).Value = 810
$ws.Cells.Item(136, 5).Value = 2426
$ws.Cells.Item(136, 7).Value = 1
$ws.Cells.Item(136, 8).Value = 57
$ws.Cells.Item(137, 1).Value = 'Sri Lanka'
$ws.Cells.Item(137, 2).Value = 3271
$ws.Cells.Item(137, 4).Value = 3021
$ws.Cells.Item(137, 5).Value = 237
$ws.Cells.Item(137, 8).Value = 13
$ws.Cells.Item(204, 1).Value = 'Santa Lucia'
$ws.Cells.Item(205, 1).Value = 'Timor Oriental'
$ws.Cells.Item(214, 1).Value = 'Montserrat'
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1
$ws.Cells.Item(215, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
